$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates from the crypto price refresh.
# Values are prefixed with a literal apostrophe so Excel stores them as
# text (matching the source inlineStr cells) instead of auto-converting
# numeric-looking strings (e.g. "223.40", "0.999") into numbers, which
# would silently drop significant trailing/formatting characters.
# The cell Style is reset to "Normal" afterwards so no stray number-format
# or quote-prefix styling is left behind on the cell.

$ws.Range("D2").Value = "'34.190.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.76%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.804.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.97%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.25%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'223.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.60%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.09%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.12%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'32.94"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +2.61%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +2.85%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0720"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.77%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0928"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'2.064.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.81%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'Chainlink"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'11.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.79%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'WrappedEther"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'1.800.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.74%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +0.44%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'34.223.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.00%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -0.79%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'68.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.30%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'247.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.12%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0789"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.86%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'11.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +5.80%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.07%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.71%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -0.23%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'159.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.27%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'16.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.64%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'7.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.26%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.63%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.26%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0529"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.62%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.04%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +1.94%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.38%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.79%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.418.65"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.89%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +3.12%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.79%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -0.91%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.946"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.75%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -3.88%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.07%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -3.82%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +4.42%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -0.22%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'108.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +4.61%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0498"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.88%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'RocketPoolETH"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'1.962.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.60%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'WEMIXToken"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'1.05"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.14%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'12.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.06%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E51").Value = "'  +2.68%  "
$ws.Range("E51").Style = "Normal"
